# perf(devops): add basic async loading readUsers IV
#
# Insert a new "YearBirth" row above the existing "Gender" row on Sheet1,
# shifting the Gender row down by one, and move the sheet's selection to A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 11, pushing the old row 11 ("Gender", ...) down to row 12.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row with the YearBirth data.
$ws.Range("A11").Value = "YearBirth"
$ws.Range("B11").Value = 1990
$ws.Range("C11").Value = 1995
$ws.Range("D11").Value = 2000

# Update the active selection to A2, matching the saved workbook state.
[void]$ws.Range("A2").Select()
